$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 35

# Column A is a text string (stored as inline/shared string in the original file,
# not a real date), so assign it directly as text.
$ws.Cells.Item($row, 1).Value = "2024-03-29 02:10:45"

$ws.Cells.Item($row, 2).Value = 31
$ws.Cells.Item($row, 3).Value = 11
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 5
$ws.Cells.Item($row, 6).Value = 6
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0.001
$ws.Cells.Item($row, 10).Value = 0.05
$ws.Cells.Item($row, 11).Value = 0.003
$ws.Cells.Item($row, 12).Value = 100
$ws.Cells.Item($row, 13).Value = 500
$ws.Cells.Item($row, 14).Value = 10
$ws.Cells.Item($row, 15).Value = 9
$ws.Cells.Item($row, 16).Value = 2
$ws.Cells.Item($row, 17).Value = 200
$ws.Cells.Item($row, 18).Value = 3
$ws.Cells.Item($row, 19).Value = 1
$ws.Cells.Item($row, 20).Value = 70
$ws.Cells.Item($row, 21).Value = 0.3548387096774194

$ws.Cells.Item($row, 22).Value = "Data/bombay1.xlsx"

$ws.Cells.Item($row, 23).Value = -32600
